$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every assignment carries a leading apostrophe so Excel treats the incoming
# text as a literal string (preserving things like leading/trailing zeros,
# "1.00" vs 1, padded percent strings, multi-dot numbers, etc.) instead of
# silently re-parsing numeric-looking text back into a number.

$ws.Range("D2").Value = "'42.304.25"
$ws.Range("E2").Value = "'  -2.79%  "
$ws.Range("D3").Value = "'2.222.60"
$ws.Range("E3").Value = "'  -2.03%  "
$ws.Range("D4").Value = "'1.01"
$ws.Range("E4").Value = "'  +0.27%  "
$ws.Range("D5").Value = "'109.40"
$ws.Range("E5").Value = "'  -8.00%  "
$ws.Range("D6").Value = "'296.43"
$ws.Range("E6").Value = "'  +11.61%  "
$ws.Range("D7").Value = "'0.624"
$ws.Range("E7").Value = "'  -3.30%  "
$ws.Range("E8").Value = "'  -0.15%  "
$ws.Range("E9").Value = "'  -2.79%  "
$ws.Range("D10").Value = "'43.96"
$ws.Range("E10").Value = "'  -7.70%  "
$ws.Range("D11").Value = "'0.0914"
$ws.Range("E11").Value = "'  -3.26%  "
$ws.Range("D12").Value = "'54.21"
$ws.Range("E12").Value = "'  +0.02%  "
$ws.Range("D13").Value = "'8.79"
$ws.Range("E13").Value = "'  -4.27%  "
$ws.Range("D14").Value = "'1.00"
$ws.Range("E14").Value = "'  +11.25%  "
$ws.Range("E15").Value = "'  -2.58%  "
$ws.Range("D16").Value = "'15.08"
$ws.Range("E16").Value = "'  -2.11%  "
$ws.Range("D17").Value = "'2.553.81"
$ws.Range("E17").Value = "'  -2.20%  "
$ws.Range("D18").Value = "'2.222.83"
$ws.Range("E18").Value = "'  -1.88%  "
$ws.Range("D19").Value = "'42.361.35"
$ws.Range("E19").Value = "'  -2.70%  "
$ws.Range("D20").Value = "'7.35"
$ws.Range("E20").Value = "'  +7.55%  "
$ws.Range("E21").Value = "'  -3.99%  "
$ws.Range("D22").Value = "'72.32"
$ws.Range("E22").Value = "'  +0.33%  "
$ws.Range("D23").Value = "'3.46"
$ws.Range("E23").Value = "'  +20.73%  "
$ws.Range("D24").Value = "'2.31"
$ws.Range("E24").Value = "'  -3.76%  "
$ws.Range("D25").Value = "'228.76"
$ws.Range("E25").Value = "'  -2.90%  "
$ws.Range("D26").Value = "'9.18"
$ws.Range("E26").Value = "'  -3.34%  "
$ws.Range("D27").Value = "'11.68"
$ws.Range("E27").Value = "'  -2.74%  "
$ws.Range("E28").Value = "'  -1.73%  "
$ws.Range("E29").Value = "'  -0.64%  "
$ws.Range("D30").Value = "'38.21"
$ws.Range("E30").Value = "'  -8.16%  "
$ws.Range("B31").Value = "'WEMIXToken"
$ws.Range("C31").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D31").Value = "'3.20"
$ws.Range("E31").Value = "'  -5.48%  "
$ws.Range("B32").Value = "'Monero"
$ws.Range("C32").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").Value = "'173.99"
$ws.Range("E32").Value = "'  +1.10%  "
$ws.Range("D33").Value = "'21.01"
$ws.Range("E33").Value = "'  -2.56%  "
$ws.Range("D34").Value = "'0.0895"
$ws.Range("E34").Value = "'  -2.05%  "
$ws.Range("D35").Value = "'5.68"
$ws.Range("E35").Value = "'  -0.32%  "
$ws.Range("D36").Value = "'5.10"
$ws.Range("E36").Value = "'  +12.00%  "
$ws.Range("D37").Value = "'4.36"
$ws.Range("E37").Value = "'  +4.38%  "
$ws.Range("D38").Value = "'0.0382"
$ws.Range("E38").Value = "'  -0.54%  "
$ws.Range("D39").Value = "'0.125"
$ws.Range("E39").Value = "'  -3.63%  "
$ws.Range("D40").Value = "'0.105"
$ws.Range("E40").Value = "'  -1.01%  "
$ws.Range("D41").Value = "'2.42"
$ws.Range("E41").Value = "'  -4.87%  "
$ws.Range("D42").Value = "'71.94"
$ws.Range("E43").Value = "'  -1.09%  "
$ws.Range("E44").Value = "'  -0.05%  "
$ws.Range("D45").Value = "'12.60"
$ws.Range("E45").Value = "'  -8.30%  "
$ws.Range("E46").Value = "'  -4.37%  "
$ws.Range("E47").Value = "'  -6.28%  "
$ws.Range("D48").Value = "'1.31"
$ws.Range("E48").Value = "'  +3.32%  "
$ws.Range("D49").Value = "'103.47"
$ws.Range("E49").Value = "'  +1.82%  "
$ws.Range("E50").Value = "'  -1.22%  "
$ws.Range("D51").Value = "'1.63"
$ws.Range("E51").Value = "'  +4.96%  "
